$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before the existing "Late" column (old column N),
# shifting Late / "heading" (Date) / Outstanding one column to the right
# (old N->O, old O->P, old P->Q).
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = 9.87

# Make "Repayment schedule" the active sheet/tab and update its cell selection
$ws.Activate()
[void]$ws.Range("S7").Select()
